{"js": "// Updates the single-column benchmark results table:\n//  - rows 1-3 (metric placeholders) become \"0M\"\n//  - row 4 count corrected 802 -> 2398\n//  - rows 5-12 (confidence-interval percentages) refreshed with new figures\n//  - the last three rows, which held a whole tab-separated results line,\n//    are collapsed back down to just their leading sample-count value\n\nconst table = context.document.body.tables.getFirst();\ntable.rows.load('items');\nawait context.sync();\n\n// 0-based row index -> new cell text\nconst updates = {\n  0: '0M',\n  1: '0M',\n  2: '0M',\n  3: '2398',\n  4: '0.00001',\n  5: '0.00067',\n  6: '0.00015',\n  7: '0.00003',\n  8: '0.00028',\n  9: '0.00033',\n  10: '0.00036',\n  11: '0.44345',\n  43: '99.88',\n  44: '0.44',\n  45: '360',\n};\n\nconst rows = table.rows.items;\nfor (const indexStr of Object.keys(updates)) {\n  const index = Number(indexStr);\n  const row = rows[index];\n  row.cells.load('items');\n}\nawait context.sync();\n\nfor (const indexStr of Object.keys(updates)) {\n  const index = Number(indexStr);\n  const row = rows[index];\n  const cell = row.cells.items[0];\n  cell.value = updates[index];\n}\nawait context.sync();\n", "ps1": "# Updates the single-column benchmark results table:\n#  - rows 1-3 (metric placeholders) become \"0M\"\n#  - row 4 count corrected 802 -> 2398\n#  - rows 5-12 (confidence-interval percentages) refreshed with new figures\n#  - the last three rows, which held a whole tab-separated results line,\n#    are collapsed back down to just their leading sample-count value\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 0-based diff row -> (1-based Word row index, new cell text)\n$updates = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"2398\"\n    5  = \"0.00001\"\n    6  = \"0.00067\"\n    7  = \"0.00015\"\n    8  = \"0.00003\"\n    9  = \"0.00028\"\n    10 = \"0.00033\"\n    11 = \"0.00036\"\n    12 = \"0.44345\"\n    44 = \"99.88\"\n    45 = \"0.44\"\n    46 = \"360\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $t.Rows.Item($rowIndex).Cells.Item(1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
